$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28, pushing existing rows 28-92 down to 29-93.
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the new data entry.
$ws.Cells.Item(28,1).Value2  = 10
$ws.Cells.Item(28,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(28,3).Value2  = "La Araucanía"
$ws.Cells.Item(28,4).Value2  = 44979
$ws.Cells.Item(28,5).Value2  = 9
$ws.Cells.Item(28,6).Value2  = "Fruta"
$ws.Cells.Item(28,7).Value2  = 100108
$ws.Cells.Item(28,8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(28,9).Value2  = 100108004
$ws.Cells.Item(28,10).Value2 = "Papaya"
$ws.Cells.Item(28,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(28,12).Value2 = "Primera"
$ws.Cells.Item(28,13).Value2 = 50
$ws.Cells.Item(28,14).Value2 = 40000
$ws.Cells.Item(28,15).Value2 = 40000
$ws.Cells.Item(28,16).Value2 = 40000
$ws.Cells.Item(28,17).Value2 = "$/caja 15 kilos granel"
$ws.Cells.Item(28,18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(28,19).Value2 = 2667
$ws.Cells.Item(28,20).Value2 = 15
